# "Actualizacion de horas. Cosas en proceso"
# Update the "TIEMPO REAL (horas)" column (E) with the latest hour tracking,
# including marking a couple of in-progress tasks as "(ongoing)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Finished tasks: replace the placeholder "(ongoing)" text with real hour counts.
$ws.Range("E5").Value = 15
$ws.Range("E6").Value = 7
$ws.Range("E7").Value = 12

# Tasks whose real-hours cell already held text stay text, just updated.
$ws.Range("E9").Value = "x"

# Newly-started tasks get an "(ongoing)" note (order matters for the shared
# string table layout, so write "media (ongoing)" before "12 (ongoing)").
$ws.Range("E12").Value = "media (ongoing)"
$ws.Range("E12").HorizontalAlignment = -4108

$ws.Range("E10").Value = "2 (ongoing)"

$ws.Range("E8").Value = "12 (ongoing)"
$ws.Range("E8").HorizontalAlignment = -4108

# Leave the view scrolled/selected roughly where the author left off.
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 6
